# Estado de Cuenta - actualizacion de base de datos:
#  - Se eliminan los periodos de mora anteriores del trabajador DIDIER LOUIS JOSEPH (CE)
#    y se agregan los nuevos periodos (2304-2310), reemplazando 2209.
#  - Se agrega una nueva trabajadora: CLARIBEL PINILLA AMAYA (CC 1007314918),
#    con periodos de mora 2304-2307.
#  - Se actualizan los totales (Valor Mora y Cantidad de Trabajadores).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 4 new blank data rows above the existing data block (row 16) ---
# This pushes the current data rows (16-23) down to (20-27) and the footer
# rows (28-29) down to (32-33), along with their merged cells, automatically.
$ws.Rows("16:19").Insert()

# --- 2. Copy the formatting of a normal data row (now row 20) onto the new rows ---
$ws.Range("B20:J20").Copy()
$ws.Range("B16:J19").PasteSpecial(-4122)

# --- 3. Fill in the new worker (CLARIBEL PINILLA AMAYA) rows 16-19 ---
$claribelPeriods = @("2307", "2306", "2305", "2304")
for ($i = 0; $i -lt 4; $i++) {
    $r = 16 + $i
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = "1007314918"
    $ws.Range("D$r").Value = "CLARIBEL PINILLA AMAYA"
    $ws.Range("E$r").Value = $claribelPeriods[$i]
    $ws.Range("F$r").Value = 72000
    $ws.Range("G$r").Value = 1800000
}

# --- 4. Update the existing worker (DIDIER LOUIS JOSEPH) rows, now 20-27 ---
$didierPeriods = @("2310", "2309", "2308", "2307", "2306", "2305", "2304", "2209")
for ($i = 0; $i -lt 8; $i++) {
    $r = 20 + $i
    $ws.Range("B$r").Value = "CE"
    $ws.Range("C$r").Value = "1090334"
    $ws.Range("D$r").Value = "DIDIER LOUIS JOSEPH"
    $ws.Range("E$r").Value = $didierPeriods[$i]
    if ($r -eq 20) {
        $ws.Range("F$r").Value = 153334
    } else {
        $ws.Range("F$r").Value = 200000
    }
    $ws.Range("G$r").Value = 5000000
}

# --- 5. Update the summary totals ---
$ws.Range("E11").Value = 1841334
$ws.Range("C13").Value = 2
